# Update column G ("K") values on Sheet1 with freshly recalculated s_vals
# (per commit: "regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals")

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$newK = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 2
    6  = 0
    7  = 1
    8  = 1
    9  = 0
    10 = 2
    11 = 0
    12 = 1
    13 = 1
    14 = 3
    15 = 0
    16 = 0
    17 = 0
    18 = 3
    19 = 2
    20 = 0
    21 = 1
    22 = 2
    23 = 2
    24 = 2
    25 = 1
    26 = 2
    27 = 1
    28 = 1
    29 = 2
    30 = 0
    31 = 2
    32 = 3
    33 = 1
    34 = 1
    35 = 5
    36 = 1
    37 = 4
    38 = 1
    39 = 2
    40 = 1
    41 = 1
    42 = 2
    43 = 6
    44 = 2
    45 = 2
    46 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
